# Add a new inferred-position observation (patrol 1, number 2) for the
# Off Chichijima Harbor sighting. This shifts all subsequent rows down by
# one (Excel renumbers rows automatically on insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 3 (Patrol 1 / #2, currently the
# "Majuro Lagoon" entry) so the new sighting becomes Patrol 1 / #2 and the
# old entries shift down to make room.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new observation.
$ws.Cells.Item(3, 1).Value2 = 1
$ws.Cells.Item(3, 2).Value2 = 2
$ws.Cells.Item(3, 3).Value2 = 800
$ws.Cells.Item(3, 4).Value2 = -9
$ws.Cells.Item(3, 5).Value2 = 16265
$ws.Cells.Item(3, 6).Value2 = 27.082999999999998
$ws.Cells.Item(3, 7).Value2 = 142.18299999999999
$ws.Cells.Item(3, 8).Value2 = "Off Chichijima Harbor"

# Update the view: select B5 (clears the old topLeftCell scroll position
# and C11 selection).
$ws.Range("B5").Select()
